$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the casing of Alexei's email address (row 3, column B):
# "OURIN@ntu.edu.sg" -> "OURIN@NTU.EDU.SG"
$ws.Range("B3").Value = "OURIN@NTU.EDU.SG"
